$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 208, shifting existing rows 208:220 down to 209:221
$ws.Rows.Item(208).Insert()

# Populate the newly inserted row 208 with the new weekly data record
$ws.Cells.Item(208, 1).Value = 11
$ws.Cells.Item(208, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(208, 3).Value = "Bíobío"
$ws.Cells.Item(208, 4).Value = 45041
$ws.Cells.Item(208, 5).Value = 8
$ws.Cells.Item(208, 6).Value = 100112032
$ws.Cells.Item(208, 7).Value = "Zapallo italiano"
$ws.Cells.Item(208, 8).Value = "Sin especificar"
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 220
$ws.Cells.Item(208, 11).Value = 4500
$ws.Cells.Item(208, 12).Value = 5000
$ws.Cells.Item(208, 13).Value = 4727
$ws.Cells.Item(208, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(208, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(208, 16).Value = 95
$ws.Cells.Item(208, 17).Value = 50
$ws.Cells.Item(208, 18).Value = "Hortaliza"
